# "error solve ifrs list"
#
# The 모나미 (Monami) IFRS "company_list" sheet had several years of
# financial data accidentally populated with figures that were too
# large by roughly three orders of magnitude (rows 2-6), and three
# forecast rows (7-9) that should never have carried any figures at
# all beyond their period label. This script corrects both issues.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2014/12): replace the mis-scaled figures with the correct ones
$ws.Range("D2").Value = 1501
$ws.Range("E2").Value = 93
$ws.Range("F2").Value = 93
$ws.Range("G2").Value = 40
$ws.Range("H2").Value = 32
$ws.Range("I2").Value = 28
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 1683
$ws.Range("L2").Value = 1109
$ws.Range("M2").Value = 574
$ws.Range("N2").Value = 537
$ws.Range("O2").Value = 37
$ws.Range("P2").Value = 139
$ws.Range("Q2").Value = 85
$ws.Range("R2").Value = 4
$ws.Range("S2").Value = -102
$ws.Range("T2").Value = 37
$ws.Range("U2").Value = 48
$ws.Range("V2").Value = 858
$ws.Range("W2").Value = 6.22
$ws.Range("X2").Value = 2.12
$ws.Range("Y2").Value = 5.34
$ws.Range("Z2").Value = 1.86
$ws.Range("AA2").Value = 193.08
$ws.Range("AB2").Value = 292.76
$ws.Range("AC2").Value = 184
$ws.Range("AD2").Value = 14.64
$ws.Range("AE2").Value = 3514
$ws.Range("AF2").Value = 0.77
$ws.Range("AG2").Value = 45
$ws.Range("AH2").Value = 1.69
$ws.Range("AI2").Value = 24.71
$ws.Range("AJ2").Value = 15283277

# Row 3 (2015/12): replace the mis-scaled figures with the correct ones
$ws.Range("D3").Value = 1429
$ws.Range("E3").Value = 97
$ws.Range("F3").Value = 97
$ws.Range("G3").Value = 65
$ws.Range("H3").Value = 51
$ws.Range("I3").Value = 47
$ws.Range("J3").Value = 4
$ws.Range("K3").Value = 1623
$ws.Range("L3").Value = 905
$ws.Range("M3").Value = 718
$ws.Range("N3").Value = 695
$ws.Range("O3").Value = 23
$ws.Range("P3").Value = 189
$ws.Range("Q3").Value = 45
$ws.Range("R3").Value = 28
$ws.Range("S3").Value = -59
$ws.Range("T3").Value = 25
$ws.Range("U3").Value = 20
$ws.Range("V3").Value = 706
$ws.Range("W3").Value = 6.76
$ws.Range("X3").Value = 3.6
$ws.Range("Y3").Value = 7.69
$ws.Range("Z3").Value = 3.12
$ws.Range("AA3").Value = 126.12
$ws.Range("AB3").Value = 272.18
$ws.Range("AC3").Value = 276
$ws.Range("AD3").Value = 16.22
$ws.Range("AE3").Value = 3678
$ws.Range("AF3").Value = 1.22
$ws.Range("AG3").Value = 60
$ws.Range("AH3").Value = 1.34
$ws.Range("AI3").Value = 23.95
$ws.Range("AJ3").Value = 18897307

# Row 4 (2016/12): replace the mis-scaled figures with the correct ones
$ws.Range("D4").Value = 1402
$ws.Range("E4").Value = 101
$ws.Range("F4").Value = 101
$ws.Range("G4").Value = 71
$ws.Range("H4").Value = 57
$ws.Range("I4").Value = 56
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1662
$ws.Range("L4").Value = 908
$ws.Range("M4").Value = 753
$ws.Range("N4").Value = 729
$ws.Range("O4").Value = 24
$ws.Range("P4").Value = 189
$ws.Range("Q4").Value = 76
$ws.Range("R4").Value = -41
$ws.Range("S4").Value = -34
$ws.Range("T4").Value = 30
$ws.Range("U4").Value = 46
$ws.Range("V4").Value = 693
$ws.Range("W4").Value = 7.2
$ws.Range("X4").Value = 4.08
$ws.Range("Y4").Value = 7.87
$ws.Range("Z4").Value = 3.48
$ws.Range("AA4").Value = 120.56
$ws.Range("AB4").Value = 294.91
$ws.Range("AC4").Value = 297
$ws.Range("AD4").Value = 14.33
$ws.Range("AE4").Value = 3907
$ws.Range("AF4").Value = 1.09
$ws.Range("AG4").Value = 70
$ws.Range("AH4").Value = 1.65
$ws.Range("AI4").Value = 23.32
$ws.Range("AJ4").Value = 18897307

# Row 5 (2017/12): replace the mis-scaled figures with the correct ones
$ws.Range("D5").Value = 1377
$ws.Range("E5").Value = 76
$ws.Range("F5").Value = 76
$ws.Range("G5").Value = 36
$ws.Range("H5").Value = 29
$ws.Range("I5").Value = 27
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 1679
$ws.Range("L5").Value = 925
$ws.Range("M5").Value = 754
$ws.Range("N5").Value = 732
$ws.Range("O5").Value = 21
$ws.Range("P5").Value = 189
$ws.Range("Q5").Value = 117
$ws.Range("R5").Value = -129
$ws.Range("S5").Value = 25
$ws.Range("T5").Value = 99
$ws.Range("U5").Value = 18
$ws.Range("V5").Value = 732
$ws.Range("W5").Value = 5.52
$ws.Range("X5").Value = 2.09
$ws.Range("Y5").Value = 3.69
$ws.Range("Z5").Value = 1.73
$ws.Range("AA5").Value = 122.65
$ws.Range("AB5").Value = 301.37
$ws.Range("AC5").Value = 143
$ws.Range("AD5").Value = 19.99
$ws.Range("AE5").Value = 3949
$ws.Range("AF5").Value = 0.72
$ws.Range("AG5").Value = 70
$ws.Range("AH5").Value = 2.46
$ws.Range("AI5").Value = 48.2
$ws.Range("AJ5").Value = 18897307

# Row 6 (2018/12): replace the mis-scaled figures with the correct ones
$ws.Range("D6").Value = 1352
$ws.Range("E6").Value = 69
$ws.Range("F6").Value = 69
$ws.Range("G6").Value = 23
$ws.Range("H6").Value = 8
$ws.Range("I6").Value = 6
$ws.Range("K6").Value = 1706
$ws.Range("L6").Value = 961
$ws.Range("M6").Value = 745
$ws.Range("N6").Value = 722
$ws.Range("P6").Value = 189
$ws.Range("Q6").Value = 74
$ws.Range("R6").Value = -75
$ws.Range("S6").Value = -5
$ws.Range("T6").Value = 103
$ws.Range("U6").Value = -30
$ws.Range("V6").Value = 746
$ws.Range("W6").Value = 5.12
$ws.Range("X6").Value = 0.5600000000000001
$ws.Range("Y6").Value = 0.84
$ws.Range("Z6").Value = 0.45
$ws.Range("AA6").Value = 129.1
$ws.Range("AB6").Value = 294.99
$ws.Range("AC6").Value = 32
$ws.Range("AD6").Value = 81.52
$ws.Range("AE6").Value = 3953
$ws.Range("AF6").Value = 0.66
$ws.Range("AG6").Value = 70
$ws.Range("AH6").Value = 2.67
$ws.Range("AI6").Value = 210.43
$ws.Range("AJ6").Value = 18897307

# Row 7: forecast period never had real figures -- clear D:AJ
# (spans every column touched by the old data, incl. already-blank ones)
$ws.Range("D7:AJ7").ClearContents()

# Row 8: forecast period never had real figures -- clear D:AJ
# (spans every column touched by the old data, incl. already-blank ones)
$ws.Range("D8:AJ8").ClearContents()

# Row 9: forecast period never had real figures -- clear D:AJ
# (spans every column touched by the old data, incl. already-blank ones)
$ws.Range("D9:AJ9").ClearContents()

